# Fruta / hortaliza, semanal
# Insert two new weekly records (rows) at the top of the data table (rows 8-9),
# pushing the existing rows 8-20 down to rows 10-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 8, shifting rows 8:20 down to 10:22.
$ws.Rows("8:9").Insert()

# --- New row 8 ---
$ws.Cells.Item(8, 1).Value  = 1
$ws.Cells.Item(8, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value  = 45125
$ws.Cells.Item(8, 5).Value  = 15
$ws.Cells.Item(8, 6).Value  = "Fruta"
$ws.Cells.Item(8, 7).Value  = 100107
$ws.Cells.Item(8, 8).Value  = "Otros"
$ws.Cells.Item(8, 9).Value  = 100107002
$ws.Cells.Item(8, 10).Value = "Chirimoya"
$ws.Cells.Item(8, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 160
$ws.Cells.Item(8, 14).Value = 14000
$ws.Cells.Item(8, 15).Value = 15000
$ws.Cells.Item(8, 16).Value = 14375
$ws.Cells.Item(8, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(8, 19).Value = 1438
$ws.Cells.Item(8, 20).Value = 10

# --- New row 9 ---
$ws.Cells.Item(9, 1).Value  = 1
$ws.Cells.Item(9, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value  = 45125
$ws.Cells.Item(9, 5).Value  = 15
$ws.Cells.Item(9, 6).Value  = "Fruta"
$ws.Cells.Item(9, 7).Value  = 100107
$ws.Cells.Item(9, 8).Value  = "Otros"
$ws.Cells.Item(9, 9).Value  = 100107002
$ws.Cells.Item(9, 10).Value = "Chirimoya"
$ws.Cells.Item(9, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 180
$ws.Cells.Item(9, 14).Value = 13000
$ws.Cells.Item(9, 15).Value = 13000
$ws.Cells.Item(9, 16).Value = 13000
$ws.Cells.Item(9, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(9, 19).Value = 1300
$ws.Cells.Item(9, 20).Value = 10
